$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12): column A labels pick up the same style used
# by the header row (row 9) for its label cell.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Update the Right/Wrong/Total numbers now that float-valued marking is handled.
$ws.Range("B10").Value = 4
$ws.Range("D10").Value = 24
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 16
$ws.Range("E12").Value = "16/112"

# --- Drop the redundant third "Student Ans / Correct Ans" block (columns G:H).
$ws.Range("G15:H21").Clear()

# --- Drop the redundant second "Student Ans / Correct Ans" block (columns D:E)
# for every question row except 16/17/18 (16 is repurposed below; 17/18 untouched).
$ws.Range("D19:E40").Clear()

# --- Row 16: D16 now carries the correct-answer text with the "correct" (green)
# style copied from B10 (which already uses that style).
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

# --- Rows 28, 30 and 39: column A now shows the matching correct answer using
# the same green "correct" style.
$ws.Range("B10").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("A39").PasteSpecial(-4122)

$ws.Range("A28").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A39").Value = "Option D"
